$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Change the text value in C8 (remove the apostrophe and double-quote)
$ws.Range("C8").Value = "< & >"

# Change the selection on the sheet from A9 to C9
$ws.Range("C9").Select()

# Change the window height setting
$excel.ActiveWindow.Height = 11760
